$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Docentes responsáveis:" section gets a new co-teacher row. Insert a
# fresh row above the existing "1033242 - Fábio Herbst Florenzano" value
# row (row 13), which pushes that row (and everything below it) down by
# one.
$ws.Rows.Item(13).Insert()

# The freshly inserted row inherits a phantom formatted cell in column A
# (copied down from the row above); clear it so row 13 has no A value,
# just like the "Docentes responsáveis:" value rows above/below it.
$ws.Range("A13").Clear()

# Fill in the new co-teacher's name in both the "current" (B) and
# "modified" (C) columns, matching the existing duplicate-value pattern
# used throughout the sheet.
$ws.Range("B13").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C13").Value = "5840963 - Daniela Camargo Vernilli"

# Match the formatting used by the sibling "Docentes responsáveis:" value
# row immediately below (now row 14): B is plain wrapped text, C is red
# wrapped text, both top-aligned.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
